$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert 6 new rows at the top of the data (row 2), pushing the existing
# glucose readings down so the most-recent readings can be added above them
# (the table is sorted newest-first).
$ws.Range("A2:A7").EntireRow.Insert()

# Grow the table definition to cover the newly inserted rows.
$lo.Resize($ws.Range("A1:C31"))

# Copy the Timestamp column's date/time number formatting down into the new
# rows (re-uses the existing style instead of fabricating a new one).
$ws.Range("A8").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New glucose readings (most recent first), matching the corrected time-of-day data.
$data = @(
    @(45987.713194444441, 10.8),
    @(45987.306944444441, 8.3000000000000007),
    @(45986.723611111112, 14.3),
    @(45986.298611111109, 8.1999999999999993),
    @(45985.887499999997, 10),
    @(45985.293055555558, 6.6)
)

$formula = '="EXEC [dbo].[NewBloodSugarReading] ''" & TEXT(Table1[[#This Row],[Timestamp]], "yyyy/mm/dd hh:mm") & "'', " & Table1[[#This Row],[Glucose Value (mmol/L)]] & ", NULL"'

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
    $ws.Cells.Item($r, 3).Formula = $formula
}

# The row-shift leaves the pre-existing calculated-column formulas (now sitting
# at rows 8..31) evaluating stale/errored until they're touched; re-assigning
# each one's own formula forces a clean re-resolution of the structured
# references against their new row positions.
for ($r = 8; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Formula = $cell.Formula
}

# Match the post-edit selection recorded in the workbook.
$ws.Range("B24").Select()
